$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new month columns: R (May-22) and S (Jun-22) ---
# Copy formatting from the existing month header (Q1) onto the new header cells,
# then set their date-serial values.
$ws.Range("Q1").Copy()
$ws.Range("R1:S1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R1").Value = 44682
$ws.Range("S1").Value = 44713

# --- Row 2 (lasgova): revise existing month, add new months ---
$ws.Range("Q2").Value = 5218
$ws.Range("R2").Value = 5254
$ws.Range("S2").Value = 5253

# --- Row 3 (lalgova): revise existing month, add new months ---
$ws.Range("Q3").Value = 14091
$ws.Range("R3").Value = 14104
$ws.Range("S3").Value = 14109

# --- Row 4 (cpgs): revise existing months, add new month ---
$ws.Range("P4").Value = 319759
$ws.Range("Q4").Value = 319650

$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R4").Value = 317251

$excel.CutCopyMode = $false

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("V20").Select()
